# Signup Page Completed 29-04-2024
# Rows 14-27 previously held the leftover numeric "SR no" sequence in
# column A (13..26). They are replaced with newly collected consumer
# sign-up email addresses in column B, matching the layout used by the
# rows above them (B2:B13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEmails = @(
    "aditya327@yopmail.com",
    "aditya549@noreply.fr",
    "adityab259@yopmail.com",
    "adityec914@yopmail.com",
    "aditygd151@yopmail.com",
    "aditycb973@yopmail.com",
    "adityaa468@yopmail.com",
    "adityge868@yopmail.com",
    "adityac769@yopmail.com",
    "aditygf387@yopmail.com",
    "adityga955@yopmail.com",
    "adityed167@yopmail.com",
    "adityde783@yopmail.com",
    "adityag363@yopmail.com"
)

$startRow = 14
for ($i = 0; $i -lt $newEmails.Length; $i++) {
    $row = $startRow + $i

    # Drop the old numeric filler that used to live in column A.
    $ws.Cells.Item($row, 1).ClearContents()

    # Write the new sign-up email address into column B.
    $ws.Cells.Item($row, 2).Value = $newEmails[$i]
}
